$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 19
$ws.Cells.Item($row, 1).Value = "Globo"
$ws.Cells.Item($row, 2).Value = "RJ TV 1"
$ws.Cells.Item($row, 3).Value = "Trânsito"
$ws.Cells.Item($row, 4).Value = "2025-04-01T12:46"
$ws.Cells.Item($row, 5).Value = "Positivo"
$ws.Cells.Item($row, 6).Value = "IMTT amplia pontos de fiscalização por videomonitoramento no município. Câmeras estão em vários pontos, entre eles Av. 28 de Março e cruzamento entre as ruas Tenente Coronel Cardoso e Barão da Lagoa Dourada. Entrevista com presidente do IMTT, Álvaro Oliveira. 140 câmeras de monitoramento instaladas, 18 cruzamentos. Monitoramento para fins de infração entre 10h da noite e 6h da manhã não estará operando o avanço de sinais. entrevista com espec. de Segurança Pública, Diogo Santana. Um homem que estava sendo procurado foi identificado pelas câmeras de monitoramento. *matéria*"

$wb.Save()
